$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = '#hier-op-spreekt'
$ws.Range('C2').Value = 'Hier op spreekt'
$ws.Range('B3').Value = '#zang-van'
$ws.Range('C3').Value = 'Zang van'
$ws.Range('B4').Value = '#seven-provintien'
$ws.Range('C4').Value = 'Seven Provintien'
$ws.Range('B5').Value = '#mars'
$ws.Range('C5').Value = 'Mars'
$ws.Range('B6').Value = '#eendragt'
$ws.Range('C6').Value = 'Eendragt'
$ws.Range('B7').Value = '#alle-de'
$ws.Range('C7').Value = 'Alle de'
$ws.Range('B8').Value = '#amstfrdam'
$ws.Range('C8').Value = 'Amstfrdam'
$ws.Range('B9').Value = '#hier-op-zegt'
$ws.Range('C9').Value = 'Hier op zegt'
$ws.Range('B10').Value = '#hercules'
$ws.Range('C10').Value = 'Hercules'
$ws.Range('B11').Value = '#hier-op-zingen-alle-de'
$ws.Range('C11').Value = 'Hier op zingen alle de'
$ws.Range('B12').Value = '#het'
$ws.Range('C12').Value = 'Het'
$ws.Range('B13').Value = '#saturnus'
$ws.Range('C13').Value = 'Saturnus'
$ws.Range('B14').Value = '#hier-op-spreekt-een-der'
$ws.Range('C14').Value = 'Hier op spreekt een der'
$ws.Range('B15').Value = '#jason'
$ws.Range('C15').Value = 'Jason'
$ws.Range('B16').Value = '#neptunus'
$ws.Range('C16').Value = 'Neptunus'
$ws.Range('B17').Value = '#lachesis'
$ws.Range('C17').Value = 'Lachesis'
$ws.Range('B18').Value = '#clotho'
$ws.Range('C18').Value = 'Clotho'
$ws.Range('B19').Value = '#pompejus'
$ws.Range('C19').Value = 'Pompejus'
$ws.Range('B20').Value = '#razerny'
$ws.Range('C20').Value = 'Razerny'
$ws.Range('B21').Value = '#triton'
$ws.Range('C21').Value = 'Triton'
$ws.Range('B22').Value = '#na-den-dans-zegt,-na-dat-hy-langzaam-naar-voren-gevoerd,-en-op-''t-strand-getreden-is'
$ws.Range('C22').Value = 'Na den dans zegt, na dat hy langzaam naar voren gevoerd, en op ''t strand getreden is'
$ws.Range('B23').Value = '#zang-van-de'
$ws.Range('C23').Value = 'Zang van de'
$ws.Range('B24').Value = '#mercurius'
$ws.Range('C24').Value = 'Mercurius'
$ws.Range('B25').Value = '#''s-lands-zeemagt'
$ws.Range('C25').Value = '''''s Lands Zeemagt'
$ws.Range('B26').Value = '#amsterdam'
$ws.Range('C26').Value = 'Amsterdam'
$ws.Range('B27').Value = '#bellona'
$ws.Range('C27').Value = 'Bellona'
$ws.Range('B28').Value = '#xerxes'
$ws.Range('C28').Value = 'Xerxes'
$ws.Range('B29').Value = '#hiero'
$ws.Range('C29').Value = 'Hiero'
$ws.Range('B30').Value = '#eölus'
$ws.Range('C30').Value = 'Eölus'
$ws.Range('B31').Value = '#een-der'
$ws.Range('C31').Value = 'Een der'
$ws.Range('B32').Value = '#atropos'
$ws.Range('C32').Value = 'Atropos'
$ws.Range('B33').Value = '#hier-na-zingt-de-zelve-met-een-der'
$ws.Range('C33').Value = 'Hier na zingt de zelve met een der'

# Clear is_prefered (D) for rows 2-32 (D33 was already empty)
$ws.Range("D2:D32").Value = $null

Write-Output "Done"
